$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1279
$ws.Range("C3").Value = 2319
$ws.Range("D3").Value = 5139
$ws.Range("E3").Value = 10500
$ws.Range("F3").Value = 14800
$ws.Range("G3").Value = 14500
$ws.Range("B4").Value = 963
$ws.Range("C4").Value = 1913
$ws.Range("D4").Value = 3855
$ws.Range("E4").Value = 10300
$ws.Range("F4").Value = 11600
$ws.Range("G4").Value = 12300
$ws.Range("B5").Value = 974
$ws.Range("C5").Value = 1889
$ws.Range("D5").Value = 4404
$ws.Range("E5").Value = 7656
$ws.Range("F5").Value = 9277
$ws.Range("G5").Value = 15900
$ws.Range("B10").Value = 15600
$ws.Range("C10").Value = 16000
$ws.Range("D10").Value = 28200
$ws.Range("E10").Value = 72600
$ws.Range("F10").Value = 134000
$ws.Range("G10").Value = 236000
$ws.Range("B11").Value = 8246
$ws.Range("C11").Value = 14500
$ws.Range("D11").Value = 27600
$ws.Range("E11").Value = 55800
$ws.Range("F11").Value = 98500
$ws.Range("G11").Value = 126000
$ws.Range("B12").Value = 8626
$ws.Range("C12").Value = 14300
$ws.Range("D12").Value = 32200
$ws.Range("E12").Value = 54900
$ws.Range("F12").Value = 95800
$ws.Range("G12").Value = 120000
$ws.Range("B17").Value = 10800
$ws.Range("C17").Value = 16400
$ws.Range("D17").Value = 19700
$ws.Range("E17").Value = 21300
$ws.Range("F17").Value = 21800
$ws.Range("G17").Value = 22200
$ws.Range("B18").Value = 7211
$ws.Range("C18").Value = 12300
$ws.Range("D18").Value = 15800
$ws.Range("E18").Value = 18300
$ws.Range("F18").Value = 19100
$ws.Range("G18").Value = 18800
$ws.Range("B19").Value = 7236
$ws.Range("C19").Value = 13200
$ws.Range("D19").Value = 13200
$ws.Range("E19").Value = 17300
$ws.Range("F19").Value = 17500
$ws.Range("G19").Value = 14400
$ws.Range("B24").Value = 358000
$ws.Range("C24").Value = 428000
$ws.Range("D24").Value = 525000
$ws.Range("E24").Value = 589000
$ws.Range("F24").Value = 618000
$ws.Range("G24").Value = 680000
$ws.Range("B25").Value = 119000
$ws.Range("C25").Value = 188000
$ws.Range("D25").Value = 278000
$ws.Range("E25").Value = 357000
$ws.Range("F25").Value = 401000
$ws.Range("G25").Value = 468000
$ws.Range("B26").Value = 143000
$ws.Range("C26").Value = 167000
$ws.Range("D26").Value = 323000
$ws.Range("E26").Value = 416000
$ws.Range("F26").Value = 384000
$ws.Range("G26").Value = 216000
$ws.Range("B31").Value = 22300
$ws.Range("C31").Value = 41800
$ws.Range("D31").Value = 42000
$ws.Range("E31").Value = 41400
$ws.Range("F31").Value = 41700
$ws.Range("G31").Value = 42000
$ws.Range("B32").Value = 3190
$ws.Range("C32").Value = 4179
$ws.Range("D32").Value = 7627
$ws.Range("E32").Value = 7827
$ws.Range("F32").Value = 8255
$ws.Range("G32").Value = 13400
$ws.Range("B33").Value = 18100
$ws.Range("C33").Value = 4496
$ws.Range("D33").Value = 5764
$ws.Range("E33").Value = 7571
$ws.Range("F33").Value = 11200
$ws.Range("G33").Value = 6858
$ws.Range("B38").Value = 110000
$ws.Range("C38").Value = 1008000
$ws.Range("D38").Value = 1074000
$ws.Range("E38").Value = 1216000
$ws.Range("F38").Value = 1226000
$ws.Range("G38").Value = 1224000
$ws.Range("B39").Value = 125000
$ws.Range("D39").Value = 226000
$ws.Range("E39").Value = 256000
$ws.Range("F39").Value = 279000
$ws.Range("G39").Value = 253000
$ws.Range("B40").Value = 136000
$ws.Range("C40").Value = 153000
$ws.Range("D40").Value = 247000
$ws.Range("E40").Value = 188000
$ws.Range("F40").Value = 616000
$ws.Range("G40").Value = 237000
$ws.Range("B45").Value = 10200
$ws.Range("C45").Value = 12300
$ws.Range("D45").Value = 13300
$ws.Range("E45").Value = 13500
$ws.Range("F45").Value = 14100
$ws.Range("G45").Value = 14600
$ws.Range("B46").Value = 7501
$ws.Range("C46").Value = 10100
$ws.Range("D46").Value = 11700
$ws.Range("E46").Value = 12900
$ws.Range("F46").Value = 13600
$ws.Range("G46").Value = 14100
$ws.Range("B47").Value = 7641
$ws.Range("C47").Value = 9416
$ws.Range("D47").Value = 11900
$ws.Range("E47").Value = 13300
$ws.Range("F47").Value = 13400
$ws.Range("G47").Value = 13100
$ws.Range("B52").Value = 302000
$ws.Range("C52").Value = 363000
$ws.Range("D52").Value = 406000
$ws.Range("E52").Value = 427000
$ws.Range("F52").Value = 447000
$ws.Range("G52").Value = 459000
$ws.Range("B53").Value = 115000
$ws.Range("C53").Value = 189000
$ws.Range("D53").Value = 257000
$ws.Range("E53").Value = 305000
$ws.Range("F53").Value = 356000
$ws.Range("G53").Value = 418000
$ws.Range("B54").Value = 132000
$ws.Range("C54").Value = 187000
$ws.Range("D54").Value = 297000
$ws.Range("E54").Value = 354000
$ws.Range("F54").Value = 366000
$ws.Range("G54").Value = 296000
